$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect("D382")

$ws.Range("D2").Value2 = 0.01691003702013555
$ws.Range("E2").Value2 = -0.02124728448047486
$ws.Range("D3").Value2 = 0.05205873288577577
$ws.Range("E3").Value2 = -0.0160733067729083
$ws.Range("D4").Value2 = 0.01499988388952956
$ws.Range("E4").Value2 = -0.0146899404880424
$ws.Range("D5").Value2 = 0.009721833552966928
$ws.Range("E5").Value2 = 0.002594210813763098
$ws.Range("D6").Value2 = 0.01573506104863808
$ws.Range("E6").Value2 = -0.00466083464792777
$ws.Range("D7").Value2 = 0.02068146199781557
$ws.Range("E7").Value2 = -0.0003016591251885359
$ws.Range("D8").Value2 = 0.004312304842097714
$ws.Range("E8").Value2 = 0.01599870717517771
$ws.Range("D9").Value2 = 0.006596882631285506
$ws.Range("E9").Value2 = 0.02028250633828321
$ws.Range("D10").Value2 = 0.0140046370637273
$ws.Range("E10").Value2 = -0.00277238702522864
$ws.Range("D11").Value2 = 0.008978492920844934
$ws.Range("E11").Value2 = 0.002956830277942135
$ws.Range("D12").Value2 = 0.01447124923463209
$ws.Range("E12").Value2 = -0.001651073197578401
$ws.Range("D13").Value2 = 0.002942002845570327
$ws.Range("E13").Value2 = -0.01894986182392411
$ws.Range("D14").Value2 = 0.006152637049043669
$ws.Range("E14").Value2 = 0.01275585879560959
$ws.Range("D15").Value2 = 0.01438344212246268
$ws.Range("E15").Value2 = 0.007761273919443479
$ws.Range("D16").Value2 = 0.01058911959852533
$ws.Range("E16").Value2 = -0.01170497814130589
$ws.Range("D17").Value2 = 0.02264169106654084
$ws.Range("E17").Value2 = -0.02842632898034569
$ws.Range("D18").Value2 = 0.008727748196902354
$ws.Range("E18").Value2 = -0.00511018843819866
$ws.Range("D19").Value2 = 0.01702781005153738
$ws.Range("E19").Value2 = 0.00972871842843781
$ws.Range("D20").Value2 = 0.01219078636680069
$ws.Range("E20").Value2 = -0.01910932055749148
$ws.Range("D21").Value2 = 0.007367036621923327
$ws.Range("E21").Value2 = -0.0003333333333334076
$ws.Range("D22").Value2 = 0.01336154786239283
$ws.Range("E22").Value2 = 0.01821974965229445
$ws.Range("D23").Value2 = 0.01941881165354728
$ws.Range("E23").Value2 = 0.01674895030683343
$ws.Range("D24").Value2 = 0.009648627774642772
$ws.Range("E24").Value2 = 0.0186515060841812
$ws.Range("D25").Value2 = 0.0208333158698062
$ws.Range("E25").Value2 = -0.002981860349540288
$ws.Range("D26").Value2 = 0.0111403531360333
$ws.Range("E26").Value2 = 0.0003127736769672484
$ws.Range("D27").Value2 = 0.01987278039673606
$ws.Range("E27").Value2 = 0.01476828125861052
$ws.Range("D28").Value2 = 0.05693060066264317
$ws.Range("E28").Value2 = -0.01999347151950381
$ws.Range("D29").Value2 = 0.01997850732771556
$ws.Range("E29").Value2 = 0.00953097567093053
$ws.Range("D30").Value2 = 0.03096279211585289
$ws.Range("E30").Value2 = -0.01714824659178604
$ws.Range("D31").Value2 = 0.01627132821895466
$ws.Range("E31").Value2 = -0.04029580440688219
$ws.Range("D32").Value2 = 0.01339838304550245
$ws.Range("E32").Value2 = -0.01204705882352941
$ws.Range("D33").Value2 = 0.02040470035174872
$ws.Range("E33").Value2 = -0.02810304449648715
$ws.Range("D34").Value2 = 0.04064473747954256
$ws.Range("E34").Value2 = -0.004310914942708854
$ws.Range("D35").Value2 = 0.0113770938534243
$ws.Range("E35").Value2 = -0.005250262513125659
$ws.Range("D36").Value2 = 0.009741478983958723
$ws.Range("E36").Value2 = -0.01697825257535279
$ws.Range("D37").Value2 = 0.01135917403461421
$ws.Range("E37").Value2 = -0.0006573181419806673
$ws.Range("D38").Value2 = 0.00730829943804583
$ws.Range("E38").Value2 = 0.002179539572265465
$ws.Range("D39").Value2 = 0.01138804485380824
$ws.Range("E39").Value2 = 0
$ws.Range("D40").Value2 = 0.01759148790766345
$ws.Range("E40").Value2 = -0.007198560287942413
$ws.Range("D41").Value2 = 0.01676279584224603
$ws.Range("E41").Value2 = 0.01143458963364163
$ws.Range("D42").Value2 = 0.03405860673954311
$ws.Range("E42").Value2 = -0.02706731752944969
$ws.Range("D43").Value2 = 0.01121210740903384
$ws.Range("E43").Value2 = -0.001715165550028264
$ws.Range("D44").Value2 = 0.02141302200831242
$ws.Range("E44").Value2 = -0.000384338217631619
$ws.Range("D45").Value2 = 0.01385978519501246
$ws.Range("E45").Value2 = -0.01463172865577722
$ws.Range("D46").Value2 = 0.007968644761197418
$ws.Range("E46").Value2 = 0.005547018477568333
$ws.Range("D47").Value2 = 0.01303666818433594
$ws.Range("E47").Value2 = 0.0136540664375715
$ws.Range("D48").Value2 = 0.009551263425774649
$ws.Range("E48").Value2 = 0.02026266416510314
$ws.Range("D49").Value2 = 0.01426756062749081
$ws.Range("E49").Value2 = 0.008624419107693626
$ws.Range("D50").Value2 = 0.008163705307430144
$ws.Range("E50").Value2 = -0.001174763319743088
$ws.Range("D51").Value2 = 0.01052198664768682
$ws.Range("E51").Value2 = 0.001438159156279983
$ws.Range("D52").Value2 = 0.008770291174151479
$ws.Range("E52").Value2 = 0.006284858070424004
$ws.Range("D53").Value2 = 0.1441332176108601
$ws.Range("E53").Value2 = 0.0001970831690973895
$ws.Range("D54").Value2 = 0.04415542909353586
$ws.Range("E54").Value2 = -0.008560112431327593
$ws.Range("E55").Value2 = -0.005166565834169301

# Update the confidential disclaimer date text (shared string used by cell A58)
$ws.Range("A58").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-24 for illustrative purposes only and are subject to change."

$ws.Protect("D382")